$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 226
$ws1.Range("F3").Value = 163
$ws1.Range("F4").Value = 150

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 226
$ws4.Range("F3").Value = 163
$ws4.Range("F4").Value = 150
$ws4.Range("F5").Value = 1
